# Case with 380 kV done: update per-line power-flow results (pl_mw) on Sheet1.
# Columns B,C,D,F,G,J,K,L,M,O for rows 2-25 get new computed values; columns
# A,E,H,I,N (indices / structurally-zero values) are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.6978504067861024
$ws.Cells.Item(2, 3).Value = 0.1930930025687729
$ws.Cells.Item(2, 4).Value = 0.1341802411594273
$ws.Cells.Item(2, 6).Value = 2.258608400063444
$ws.Cells.Item(2, 7).Value = 0.002533610793596354
$ws.Cells.Item(2, 10).Value = 0.2647243273998949
$ws.Cells.Item(2, 11).Value = 0.3129193531068495
$ws.Cells.Item(2, 12).Value = 0.2973886175894407
$ws.Cells.Item(2, 13).Value = 0.2204887186071431
$ws.Cells.Item(2, 15).Value = 5.997299208462749
$ws.Cells.Item(3, 2).Value = 0.6661859499192815
$ws.Cells.Item(3, 3).Value = 0.1928607087091514
$ws.Cells.Item(3, 4).Value = 0.1328966867213168
$ws.Cells.Item(3, 6).Value = 2.268062712592823
$ws.Cells.Item(3, 7).Value = 0.002536200519294159
$ws.Cells.Item(3, 10).Value = 0.2662980842350073
$ws.Cells.Item(3, 11).Value = 0.2839817536946043
$ws.Cells.Item(3, 12).Value = 0.2954499497816911
$ws.Cells.Item(3, 13).Value = 0.2147544219842494
$ws.Cells.Item(3, 15).Value = 6.02915790192975
$ws.Cells.Item(4, 2).Value = 0.6470008946372729
$ws.Cells.Item(4, 3).Value = 0.1927276880475262
$ws.Cells.Item(4, 4).Value = 0.1321536003543429
$ws.Cells.Item(4, 6).Value = 2.27476649796732
$ws.Cells.Item(4, 7).Value = 0.002537876629307612
$ws.Cells.Item(4, 10).Value = 0.2673294873770597
$ws.Cells.Item(4, 11).Value = 0.2662698578183011
$ws.Cells.Item(4, 12).Value = 0.2943630801324915
$ws.Cells.Item(4, 13).Value = 0.2113186052211873
$ws.Cells.Item(4, 15).Value = 6.051137240658676
$ws.Cells.Item(5, 2).Value = 0.6392480901513125
$ws.Cells.Item(5, 3).Value = 0.1926759165077065
$ws.Cells.Item(5, 4).Value = 0.1318621667292135
$ws.Cells.Item(5, 6).Value = 2.277724616218279
$ws.Cells.Item(5, 7).Value = 0.002538581348532924
$ws.Cells.Item(5, 10).Value = 0.2677661908497182
$ws.Cells.Item(5, 11).Value = 0.2590666248657243
$ws.Cells.Item(5, 12).Value = 0.2939462945423585
$ws.Cells.Item(5, 13).Value = 0.2099400122595156
$ws.Cells.Item(5, 15).Value = 6.060702532090488
$ws.Cells.Item(6, 2).Value = 0.63796470289833
$ws.Cells.Item(6, 3).Value = 0.192667467550077
$ws.Cells.Item(6, 4).Value = 0.1318144634696381
$ws.Cells.Item(6, 6).Value = 2.278229483404147
$ws.Cells.Item(6, 7).Value = 0.002538699678624444
$ws.Cells.Item(6, 10).Value = 0.2678396963072309
$ws.Cells.Item(6, 11).Value = 0.257871423076935
$ws.Cells.Item(6, 12).Value = 0.2938786683478227
$ws.Cells.Item(6, 13).Value = 0.2097124024519381
$ws.Cells.Item(6, 15).Value = 6.062327610432845
$ws.Cells.Item(7, 2).Value = 0.646896072535128
$ws.Cells.Item(7, 3).Value = 0.1927269799519742
$ws.Cells.Item(7, 4).Value = 0.132149623821654
$ws.Cells.Item(7, 6).Value = 2.274805475675976
$ws.Cells.Item(7, 7).Value = 0.002537886045514188
$ws.Cells.Item(7, 10).Value = 0.2673353104819682
$ws.Cells.Item(7, 11).Value = 0.2661726531574118
$ws.Cells.Item(7, 12).Value = 0.2943573533107866
$ws.Cells.Item(7, 13).Value = 0.2112999256728223
$ws.Cells.Item(7, 15).Value = 6.051263776923093
$ws.Cells.Item(8, 2).Value = 0.6868795143652164
$ws.Cells.Item(8, 3).Value = 0.1930109247348923
$ws.Cells.Item(8, 4).Value = 0.1337283622343435
$ws.Cells.Item(8, 6).Value = 2.261681855790592
$ws.Cells.Item(8, 7).Value = 0.002534485918384646
$ws.Cells.Item(8, 10).Value = 0.2652534607612029
$ws.Cells.Item(8, 11).Value = 0.3029303257124525
$ws.Cells.Item(8, 12).Value = 0.2966987473922629
$ws.Cells.Item(8, 13).Value = 0.218493959716902
$ws.Cells.Item(8, 15).Value = 6.00778254058153
$ws.Cells.Item(9, 2).Value = 0.7673035586720687
$ws.Cells.Item(9, 3).Value = 0.1936431910138907
$ws.Cells.Item(9, 4).Value = 0.1371791485749014
$ws.Cells.Item(9, 6).Value = 2.243067553785217
$ws.Cells.Item(9, 7).Value = 0.002528497812785841
$ws.Cells.Item(9, 10).Value = 0.2616864059621502
$ws.Cells.Item(9, 11).Value = 0.3754395301587294
$ws.Cells.Item(9, 12).Value = 0.3021072933550855
$ws.Cells.Item(9, 13).Value = 0.2332711337137887
$ws.Cells.Item(9, 15).Value = 5.941680815107475
$ws.Cells.Item(10, 2).Value = 0.8275962896017859
$ws.Cells.Item(10, 3).Value = 0.1941527373338445
$ws.Cells.Item(10, 4).Value = 0.1399280286814388
$ws.Cells.Item(10, 6).Value = 2.233720157346781
$ws.Cells.Item(10, 7).Value = 0.0025245085484639
$ws.Cells.Item(10, 10).Value = 0.2593782885805815
$ws.Cells.Item(10, 11).Value = 0.4289565102831716
$ws.Cells.Item(10, 12).Value = 0.3065744053036212
$ws.Cells.Item(10, 13).Value = 0.2445304533722492
$ws.Cells.Item(10, 15).Value = 5.904774293181305
$ws.Cells.Item(11, 2).Value = 0.8552819982160997
$ws.Cells.Item(11, 3).Value = 0.1943941200324346
$ws.Cells.Item(11, 4).Value = 0.1412243840161977
$ws.Cells.Item(11, 6).Value = 2.230405114013891
$ws.Cells.Item(11, 7).Value = 0.002522781934863844
$ws.Cells.Item(11, 10).Value = 0.2583957981860969
$ws.Cells.Item(11, 11).Value = 0.4533527750027417
$ws.Cells.Item(11, 12).Value = 0.3087128176907754
$ws.Cells.Item(11, 13).Value = 0.2497388846022659
$ws.Cells.Item(11, 15).Value = 5.890511038369226
$ws.Cells.Item(12, 2).Value = 0.8658024516745684
$ws.Cells.Item(12, 3).Value = 0.1944868863180602
$ws.Cells.Item(12, 4).Value = 0.1417218224920873
$ws.Cells.Item(12, 6).Value = 2.229284320853182
$ws.Cells.Item(12, 7).Value = 0.002522140717982325
$ws.Cells.Item(12, 10).Value = 0.2580334328288458
$ws.Cells.Item(12, 11).Value = 0.4625979773294944
$ws.Cells.Item(12, 12).Value = 0.3095377722176949
$ws.Cells.Item(12, 13).Value = 0.2517234990440897
$ws.Cells.Item(12, 15).Value = 5.885472669250333
$ws.Cells.Item(13, 2).Value = 0.8635350719565906
$ws.Cells.Item(13, 3).Value = 0.1944668472534445
$ws.Cells.Item(13, 4).Value = 0.1416144004424922
$ws.Cells.Item(13, 6).Value = 2.229519723384129
$ws.Cells.Item(13, 7).Value = 0.00252227825547783
$ws.Cells.Item(13, 10).Value = 0.2581110444063803
$ws.Cells.Item(13, 11).Value = 0.4606065599146802
$ws.Cells.Item(13, 12).Value = 0.3093594297983486
$ws.Cells.Item(13, 13).Value = 0.251295532323816
$ws.Cells.Item(13, 15).Value = 5.886541641857008
$ws.Cells.Item(14, 2).Value = 0.8561467949773487
$ws.Cells.Item(14, 3).Value = 0.1944017248222352
$ws.Cells.Item(14, 4).Value = 0.1412651779152156
$ws.Cells.Item(14, 6).Value = 2.230310210510112
$ws.Cells.Item(14, 7).Value = 0.002522728929080321
$ws.Cells.Item(14, 10).Value = 0.2583657922433282
$ws.Cells.Item(14, 11).Value = 0.4541132485234982
$ws.Cells.Item(14, 12).Value = 0.3087803834650771
$ws.Cells.Item(14, 13).Value = 0.2499019144575385
$ws.Cells.Item(14, 15).Value = 5.890089259498694
$ws.Cells.Item(15, 2).Value = 0.8516259927897352
$ws.Cells.Item(15, 3).Value = 0.1943620119989475
$ws.Cells.Item(15, 4).Value = 0.141052118543314
$ws.Cells.Item(15, 6).Value = 2.230811921577555
$ws.Cells.Item(15, 7).Value = 0.002523006620862532
$ws.Cells.Item(15, 10).Value = 0.2585230929472999
$ws.Cells.Item(15, 11).Value = 0.4501367876376889
$ws.Cells.Item(15, 12).Value = 0.3084276752764055
$ws.Cells.Item(15, 13).Value = 0.2490498800780756
$ws.Cells.Item(15, 15).Value = 5.892309515887064
$ws.Cells.Item(16, 2).Value = 0.8257921036111782
$ws.Cells.Item(16, 3).Value = 0.1941371537221315
$ws.Cells.Item(16, 4).Value = 0.1398442268159243
$ws.Cells.Item(16, 6).Value = 2.233955642211328
$ws.Cells.Item(16, 7).Value = 0.002524623154993526
$ws.Cells.Item(16, 10).Value = 0.2594438527302927
$ws.Cells.Item(16, 11).Value = 0.427363144228849
$ws.Cells.Item(16, 12).Value = 0.3064367853647525
$ws.Cells.Item(16, 13).Value = 0.244191798505085
$ws.Cells.Item(16, 15).Value = 5.905757221806596
$ws.Cells.Item(17, 2).Value = 0.8100095224981771
$ws.Cells.Item(17, 3).Value = 0.1940016526187378
$ws.Cells.Item(17, 4).Value = 0.1391149315210072
$ws.Cells.Item(17, 6).Value = 2.236124081129688
$ws.Cells.Item(17, 7).Value = 0.002525637375330617
$ws.Cells.Item(17, 10).Value = 0.2600259777912193
$ws.Cells.Item(17, 11).Value = 0.4134049997290674
$ws.Cells.Item(17, 12).Value = 0.305242592593217
$ws.Cells.Item(17, 13).Value = 0.2412335765395852
$ws.Cells.Item(17, 15).Value = 5.914653579016374
$ws.Cells.Item(18, 2).Value = 0.8009561455891117
$ws.Cells.Item(18, 3).Value = 0.1939246194989934
$ws.Cells.Item(18, 4).Value = 0.1386997818051299
$ws.Cells.Item(18, 6).Value = 2.237459530301663
$ws.Cells.Item(18, 7).Value = 0.002526229026149446
$ws.Cells.Item(18, 10).Value = 0.2603671539073495
$ws.Cells.Item(18, 11).Value = 0.4053814849398236
$ws.Cells.Item(18, 12).Value = 0.3045657360971177
$ws.Cells.Item(18, 13).Value = 0.2395402369389146
$ws.Cells.Item(18, 15).Value = 5.920008282394207
$ws.Cells.Item(19, 2).Value = 0.7978950293199603
$ws.Cells.Item(19, 3).Value = 0.1938986931669717
$ws.Cells.Item(19, 4).Value = 0.1385599633242833
$ws.Cells.Item(19, 6).Value = 2.237926849010393
$ws.Cells.Item(19, 7).Value = 0.002526430775939874
$ws.Cells.Item(19, 10).Value = 0.2604837621048333
$ws.Cells.Item(19, 11).Value = 0.4026657068977499
$ws.Cells.Item(19, 12).Value = 0.3043382867234499
$ws.Cells.Item(19, 13).Value = 0.2389683057002614
$ws.Cells.Item(19, 15).Value = 5.92186213896494
$ws.Cells.Item(20, 2).Value = 0.8116870902244955
$ws.Cells.Item(20, 3).Value = 0.1940159835672688
$ws.Cells.Item(20, 4).Value = 0.1391921194125558
$ws.Cells.Item(20, 6).Value = 2.235884118188892
$ws.Cells.Item(20, 7).Value = 0.002525528551474577
$ws.Cells.Item(20, 10).Value = 0.2599633522122922
$ws.Cells.Item(20, 11).Value = 0.4148903705576856
$ws.Cells.Item(20, 12).Value = 0.305368681027943
$ws.Cells.Item(20, 13).Value = 0.24154764173052
$ws.Cells.Item(20, 15).Value = 5.913681942688015
$ws.Cells.Item(21, 2).Value = 0.8583159263011169
$ws.Cells.Item(21, 3).Value = 0.1944208161016618
$ws.Cells.Item(21, 4).Value = 0.1413675761069015
$ws.Cells.Item(21, 6).Value = 2.230074375581893
$ws.Cells.Item(21, 7).Value = 0.002522596213477824
$ws.Cells.Item(21, 10).Value = 0.2582907040215119
$ws.Cells.Item(21, 11).Value = 0.4560203089867798
$ws.Cells.Item(21, 12).Value = 0.3089500521255815
$ws.Cells.Item(21, 13).Value = 0.2503109212934973
$ws.Cells.Item(21, 15).Value = 5.889037394487644
$ws.Cells.Item(22, 2).Value = 0.8890028744148708
$ws.Cells.Item(22, 3).Value = 0.1946933155214978
$ws.Cells.Item(22, 4).Value = 0.1428274367735014
$ws.Cells.Item(22, 6).Value = 2.227061484965176
$ws.Cells.Item(22, 7).Value = 0.002520753259559996
$ws.Cells.Item(22, 10).Value = 0.2572539584258386
$ws.Cells.Item(22, 11).Value = 0.4829408996277209
$ws.Cells.Item(22, 12).Value = 0.3113791430375841
$ws.Cells.Item(22, 13).Value = 0.2561098259542689
$ws.Cells.Item(22, 15).Value = 5.875045351783911
$ws.Cells.Item(23, 2).Value = 0.8726054636578056
$ws.Cells.Item(23, 3).Value = 0.1945471591070174
$ws.Cells.Item(23, 4).Value = 0.1420448176999969
$ws.Cells.Item(23, 6).Value = 2.228597846732001
$ws.Cells.Item(23, 7).Value = 0.002521730172470172
$ws.Cells.Item(23, 10).Value = 0.2578021331292959
$ws.Cells.Item(23, 11).Value = 0.468569392089762
$ws.Cells.Item(23, 12).Value = 0.3100746311653211
$ws.Cells.Item(23, 13).Value = 0.2530083384575335
$ws.Cells.Item(23, 15).Value = 5.882319807377968
$ws.Cells.Item(24, 2).Value = 0.8109285987633541
$ws.Cells.Item(24, 3).Value = 0.1940095018386288
$ws.Cells.Item(24, 4).Value = 0.1391572098925593
$ws.Cells.Item(24, 6).Value = 2.235992328857833
$ws.Cells.Item(24, 7).Value = 0.002525577724123328
$ws.Cells.Item(24, 10).Value = 0.2599916449732724
$ws.Cells.Item(24, 11).Value = 0.4142188306065862
$ws.Cells.Item(24, 12).Value = 0.3053116462221794
$ws.Cells.Item(24, 13).Value = 0.2414056298541709
$ws.Cells.Item(24, 15).Value = 5.914120471635641
$ws.Cells.Item(25, 2).Value = 0.7453333384009397
$ws.Cells.Item(25, 3).Value = 0.1934641704523727
$ws.Cells.Item(25, 4).Value = 0.1362079124537701
$ws.Cells.Item(25, 6).Value = 2.247342166294665
$ws.Cells.Item(25, 7).Value = 0.002530045427112299
$ws.Cells.Item(25, 10).Value = 0.2625963760524979
$ws.Cells.Item(25, 11).Value = 0.3557797766669069
$ws.Cells.Item(25, 12).Value = 0.3005571470691281
$ws.Cells.Item(25, 13).Value = 0.2292024152742975
$ws.Cells.Item(25, 15).Value = 5.957513942933616
